$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '22.360.45'
$ws.Range('E2').Value = '  -4.36%  '
$ws.Range('D3').Value = '1.568.20'
$ws.Range('E3').Value = '  -3.81%  '
$ws.Range('E4').Value = '  +0.17%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '1.002'
$ws.Range('E5').Value = '  +0.11%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '289.35'
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.3679'
$ws.Range('E7').Value = '  -2.16%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '49.33'
$ws.Range('E8').Value = '  -2.13%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.3375'
$ws.Range('E9').Value = '  -3.18%  '
$ws.Range('E10').Value = '  -3.42%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07617'
$ws.Range('E11').Value = '  -5.24%  '
$ws.Range('E12').Value = '  +0.14%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '21.25'
$ws.Range('E13').Value = '  -2.84%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '6.061'
$ws.Range('E14').Value = '  -4.05%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '6.904'
$ws.Range('D16').Value = '1.569.10'
$ws.Range('E16').Value = '  -3.73%  '
$ws.Range('E17').Value = '  -5.17%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '89.47'
$ws.Range('E18').Value = '  -5.49%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.06747'
$ws.Range('E19').Value = '  -2.81%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '1.001'
$ws.Range('E20').Value = '  +0.03%  '
$ws.Range('E21').Value = '  -6.11%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.5324'
$ws.Range('E22').Value = '  -6.01%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '16.48'
$ws.Range('E23').Value = '  -4.78%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '11.97'
$ws.Range('E24').Value = '  -3.32%  '
$ws.Range('D25').Value = '22.379.22'
$ws.Range('E25').Value = '  -4.29%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.375'
$ws.Range('E26').Value = '  -2.00%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '2.901'
$ws.Range('E27').Value = '  -2.66%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '19.99'
$ws.Range('E28').Value = '  -3.69%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '146.01'
$ws.Range('E29').Value = '  -4.11%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '4.967'
$ws.Range('E30').Value = '  -4.18%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '125.23'
$ws.Range('E31').Value = '  -4.92%  '
$ws.Range('D32').Value = '1.744.61'
$ws.Range('E32').Value = '  -3.65%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.035'
$ws.Range('E33').Value = '  +6.40%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '6.253'
$ws.Range('E34').Value = '  -7.46%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '2.004'
$ws.Range('E35').Value = '  -5.17%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '10.21'
$ws.Range('E36').Value = '  -8.77%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.08479'
$ws.Range('E37').Value = '  -2.75%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.02534'
$ws.Range('E38').Value = '  -5.08%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.2321'
$ws.Range('E39').Value = '  -4.32%  '
$ws.Range('B40').Value = 'InternetComputer(DFINITY)'
$ws.Range('C40').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '5.533'
$ws.Range('E40').Value = '  -5.47%  '
$ws.Range('B41').Value = 'Hedera'
$ws.Range('C41').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.06521'
$ws.Range('E41').Value = '  -2.50%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.268'
$ws.Range('E42').Value = '  -1.21%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '11.68'
$ws.Range('E43').Value = '  -8.21%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.6349'
$ws.Range('E44').Value = '  -6.60%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '14.35'
$ws.Range('E45').Value = '  -6.43%  '
$ws.Range('E46').Value = '  +0.11%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.5975'
$ws.Range('E47').Value = '  -5.38%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '3.746'
$ws.Range('E48').Value = '  -3.70%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.105'
$ws.Range('E49').Value = '  -5.62%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.255'
$ws.Range('E50').Value = '  +3.37%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '124.29'
$ws.Range('E51').Value = '  -2.01%  '
